# edit.ps1 - apply the syllabus.docx changes described in the commit diff:
#   1. Change the date text "1/1/23" -> "2023-01-01"
#   2. Add a new "Abstract Title" paragraph style
#   3. Change the "Abstract" style's space-before from 15pt (300) to 5pt (100)
#   4. Add a new "Footnote Block Text" paragraph style

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the date on the title page.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("1/1/23", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2023-01-01", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Add the "Abstract Title" style (based on Normal, followed by Abstract).
# ---------------------------------------------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060   # BGR-packed value of RGB 34 5A 8A

# ---------------------------------------------------------------------------
# 3. Tighten the space above the "Abstract" style's paragraphs.
# ---------------------------------------------------------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5
$abstract.ParagraphFormat.SpaceAfter = 15

# ---------------------------------------------------------------------------
# 4. Add the "Footnote Block Text" style (based on Footnote Text).
# ---------------------------------------------------------------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = $d.Styles("Footnote Text")
$footnoteBlockText.NextParagraphStyle = $d.Styles("Footnote Text")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24
